$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 3 (R) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 367
$wsOff.Range("C3").Value = 261
$wsOff.Range("D3").Value = 94
$wsOff.Range("E3").Value = 46

# Sheet "DEF" - row 3 (R) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 391
$wsDef.Range("C3").Value = 268
$wsDef.Range("D3").Value = 93
$wsDef.Range("E3").Value = 47
$wsDef.Range("G3").Value = 4
